# Update odds/stats values for the FlashScore "Jogos da Semana" sheet and
# drop the last match (row 11 - "UTA Arad vs Univ. Craiova") which was
# removed from the source feed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Newcastle - West Ham) ---
$ws.Range("G2").Value  = 1.48
$ws.Range("H2").Value  = 4.5
$ws.Range("I2").Value  = 6.25
$ws.Range("J2").Value  = 1.92
$ws.Range("AD2").Value = 9.5
$ws.Range("AE2").Value = 17
$ws.Range("AH2").Value = 23
$ws.Range("AQ2").Value = 21

# --- Row 3 ---
$ws.Range("M3").Value  = 1.08
$ws.Range("N3").Value  = 8
$ws.Range("Q3").Value  = 2.15
$ws.Range("R3").Value  = 1.62
$ws.Range("AT3").Value = 2.62

# --- Row 6 ---
$ws.Range("G6").Value  = 1.42
$ws.Range("H6").Value  = 5
$ws.Range("I6").Value  = 6
$ws.Range("J6").Value  = 1.91
$ws.Range("K6").Value  = 2.63
$ws.Range("S6").Value  = 1.22
$ws.Range("U6").Value  = 1.58
$ws.Range("V6").Value  = 2.2
$ws.Range("AB6").Value = 21
$ws.Range("AD6").Value = 10
$ws.Range("AG6").Value = 151
$ws.Range("AH6").Value = 21
$ws.Range("AI6").Value = 34
$ws.Range("AM6").Value = 34
$ws.Range("AO6").Value = 7
$ws.Range("AP6").Value = 15
$ws.Range("AX6").Value = 8

# --- Row 7 ---
$ws.Range("S7").Value = 1.29
$ws.Range("U7").Value = 1.47

# --- Row 8 ---
$ws.Range("G8").Value  = 1.36
$ws.Range("I8").Value  = 6
$ws.Range("S8").Value  = 1.14
$ws.Range("U8").Value  = 1.37
$ws.Range("X8").Value  = 13
$ws.Range("AE8").Value = 15
$ws.Range("AK8").Value = 67
$ws.Range("AQ8").Value = 15
$ws.Range("AW8").Value = 151
$ws.Range("AX8").Value = 9

# --- Row 9 ---
$ws.Range("G9").Value  = 4.33
$ws.Range("I9").Value  = 1.65
$ws.Range("Q9").Value  = 1.5
$ws.Range("S9").Value  = 1.25
$ws.Range("T9").Value  = 3.75
$ws.Range("U9").Value  = 1.54
$ws.Range("V9").Value  = 2.25
$ws.Range("Y9").Value  = 15
$ws.Range("Z9").Value  = 51
$ws.Range("AK9").Value = 13
$ws.Range("AO9").Value = 23
$ws.Range("AP9").Value = 26
$ws.Range("AR9").Value = 81
$ws.Range("AT9").Value = 3.75

# --- Remove the last match row entirely (row 11) ---
$ws.Rows.Item(11).Delete()
